$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update loading_percent data for rows 2-25 (case with 380 kV done)
$data = New-Object 'object[,]' 24,12

$data[0,0] = 16.69781362728563
$data[0,1] = 14.19332454900518
$data[0,2] = 6.951549914022692
$data[0,3] = 0
$data[0,4] = 46.29031135036597
$data[0,5] = 3.70747133614231
$data[0,6] = 0
$data[0,7] = 30.09249396601462
$data[0,8] = 10.72921214413846
$data[0,9] = 0
$data[0,10] = 0
$data[0,11] = 19.52158227580475

$data[1,0] = 16.32793834108545
$data[1,1] = 13.79814439600413
$data[1,2] = 6.961493729702431
$data[1,3] = 0
$data[1,4] = 46.01280776356245
$data[1,5] = 3.711618159176107
$data[1,6] = 0
$data[1,7] = 30.00784621329311
$data[1,8] = 10.74991373997399
$data[1,9] = 0
$data[1,10] = 0
$data[1,11] = 19.43460482876289

$data[2,0] = 16.10222227990568
$data[2,1] = 13.5547133967135
$data[2,2] = 6.968100583810916
$data[2,3] = 0
$data[2,4] = 45.85352756971746
$data[2,5] = 3.714294277338145
$data[2,6] = 0
$data[2,7] = 29.96236623978795
$data[2,8] = 10.76389735756344
$data[2,9] = 0
$data[2,10] = 0
$data[2,11] = 19.38601484201444

$data[3,0] = 16.0107420154122
$data[3,1] = 13.45548810222221
$data[3,2] = 6.970919492493937
$data[3,3] = 0
$data[3,4] = 45.79145097569816
$data[3,5] = 3.715417626598348
$data[3,6] = 0
$data[3,7] = 29.94546807549523
$data[3,8] = 10.7699157711012
$data[3,9] = 0
$data[3,10] = 0
$data[3,11] = 19.36743982114913

$data[4,0] = 15.9955865192064
$data[4,1] = 13.43901550240122
$data[4,2] = 6.971395229043694
$data[4,3] = 0
$data[4,4] = 45.78131522924761
$data[4,5] = 3.715606143361751
$data[4,6] = 0
$data[4,7] = 29.94276094815906
$data[4,8] = 10.77093444773968
$data[4,9] = 0
$data[4,10] = 0
$data[4,11] = 19.36442989761392

$data[5,0] = 16.10098631491047
$data[5,1] = 13.55337507731576
$data[5,2] = 6.968138087470542
$data[5,3] = 0
$data[5,4] = 45.85267887253006
$data[5,5] = 3.71430929420322
$data[5,6] = 0
$data[5,7] = 29.96213172248471
$data[5,8] = 10.76397722852572
$data[5,9] = 0
$data[5,10] = 0
$data[5,11] = 19.38575935090427

$data[6,0] = 16.5700824151907
$data[6,1] = 14.05733160174982
$data[6,2] = 6.954874795001587
$data[6,3] = 0
$data[6,4] = 46.19234699608722
$data[6,5] = 3.708874268758787
$data[6,6] = 0
$data[6,7] = 30.06196147107842
$data[6,8] = 10.73608585442338
$data[6,9] = 0
$data[6,10] = 0
$data[6,11] = 19.49060278381311

$data[7,0] = 17.49451482725167
$data[7,1] = 15.03214775083697
$data[7,2] = 6.932821225749021
$data[7,3] = 0
$data[7,4] = 46.9446722764034
$data[7,5] = 3.699241291082764
$data[7,6] = 0
$data[7,7] = 30.30908020832248
$data[7,8] = 10.69149299051394
$data[7,9] = 0
$data[7,10] = 0
$data[7,11] = 19.73366071170085

$data[8,0] = 18.16813390252603
$data[8,1] = 15.73121953894134
$data[8,2] = 6.919000633752481
$data[8,3] = 0
$data[8,4] = 47.54712958681639
$data[8,5] = 3.692780431406797
$data[8,6] = 0
$data[8,7] = 30.52160508920011
$data[8,8] = 10.66489395711661
$data[8,9] = 0
$data[8,10] = 0
$data[8,11] = 19.93397675907671

$data[9,0] = 18.47173968167431
$data[9,1] = 16.04385168594498
$data[9,2] = 6.913224446298408
$data[9,3] = 0
$data[9,4] = 47.83130740353661
$data[9,5] = 3.68997327793588
$data[9,6] = 0
$data[9,7] = 30.62490384058389
$data[9,8] = 10.65413268364646
$data[9,9] = 0
$data[9,10] = 0
$data[9,11] = 20.02956354984761

$data[10,0] = 18.58617514096417
$data[10,1] = 16.16133954120925
$data[10,2] = 6.91111013141526
$data[10,3] = 0
$data[10,4] = 47.94030900091504
$data[10,5] = 3.688929113097758
$data[10,6] = 0
$data[10,7] = 30.66496024191033
$data[10,8] = 10.65025025831776
$data[10,9] = 0
$data[10,10] = 0
$data[10,11] = 20.06637728381386

$data[11,0] = 18.56155507536874
$data[11,1] = 16.13607830027171
$data[11,2] = 6.91156224743821
$data[11,3] = 0
$data[11,4] = 47.9167728225606
$data[11,5] = 3.689153156642166
$data[11,6] = 0
$data[11,7] = 30.65629181622441
$data[11,8] = 10.65107783867977
$data[11,9] = 0
$data[11,10] = 0
$data[11,11] = 20.05842174463148

$data[12,0] = 18.48116574454263
$data[12,1] = 16.05353618245618
$data[12,2] = 6.913049039885721
$data[12,3] = 0
$data[12,4] = 47.84024756453624
$data[12,5] = 3.689886996947025
$data[12,6] = 0
$data[12,7] = 30.62818055551954
$data[12,8] = 10.65380941279789
$data[12,9] = 0
$data[12,10] = 0
$data[12,11] = 20.03258001495712

$data[13,0] = 18.43185186079784
$data[13,1] = 16.00285611933515
$data[13,2] = 6.913969236649519
$data[13,3] = 0
$data[13,4] = 47.7935526179942
$data[13,5] = 3.690338945499991
$data[13,6] = 0
$data[13,7] = 30.61108353417568
$data[13,8] = 10.6555076711666
$data[13,9] = 0
$data[13,10] = 0
$data[13,11] = 20.01683082557988

$data[14,0] = 18.14822462399625
$data[14,1] = 15.71066905280952
$data[14,2] = 6.919388361670744
$data[14,3] = 0
$data[14,4] = 47.52875597664731
$data[14,5] = 3.692966529251249
$data[14,6] = 0
$data[14,7] = 30.51498645082701
$data[14,8] = 10.66562418162743
$data[14,9] = 0
$data[14,10] = 0
$data[14,11] = 19.92781765246261

$data[15,0] = 17.97341150713261
$data[15,1] = 15.52995104141065
$data[15,2] = 6.922843335480862
$data[15,3] = 0
$data[15,4] = 47.36885764081806
$data[15,5] = 3.694612165036387
$data[15,6] = 0
$data[15,7] = 30.45772116674279
$data[15,8] = 10.67217331169494
$data[15,9] = 0
$data[15,10] = 0
$data[15,11] = 19.87433668401652

$data[16,0] = 17.87260462120409
$data[16,1] = 15.42550730314848
$data[16,2] = 6.924878657573257
$data[16,3] = 0
$data[16,4] = 47.27784471840512
$data[16,5] = 3.69557111611743
$data[16,6] = 0
$data[16,7] = 30.42540766622824
$data[16,8] = 10.6760662016543
$data[16,9] = 0
$data[16,10] = 0
$data[16,11] = 19.84399737645663

$data[17,0] = 17.8384326614882
$data[17,1] = 15.39006266907595
$data[17,2] = 6.925576060175652
$data[17,3] = 0
$data[17,4] = 47.2471954987111
$data[17,5] = 3.695897938346014
$data[17,6] = 0
$data[17,7] = 30.4145743940161
$data[17,8] = 10.67740590465071
$data[17,9] = 0
$data[17,10] = 0
$data[17,11] = 19.83379813571664

$data[18,0] = 17.99204841395089
$data[18,1] = 15.54924142648318
$data[18,2] = 6.922470571512818
$data[18,3] = 0
$data[18,4] = 47.38578057319027
$data[18,5] = 3.694435699290374
$data[18,6] = 0
$data[18,7] = 30.46375266055347
$data[18,8] = 10.67146310381704
$data[18,9] = 0
$data[18,10] = 0
$data[18,11] = 19.87998637178235

$data[19,0] = 18.50479348788257
$data[19,1] = 16.07780615147862
$data[19,2] = 6.912610355288562
$data[19,3] = 0
$data[19,4] = 47.86268770827131
$data[19,5] = 3.689670939872983
$data[19,6] = 0
$data[19,7] = 30.63641212723115
$data[19,8] = 10.65300185424882
$data[19,9] = 0
$data[19,10] = 0
$data[19,11] = 20.04015381073635

$data[20,0] = 18.83673969162877
$data[20,1] = 16.41796242823128
$data[20,2] = 6.906591436290191
$data[20,3] = 0
$data[20,4] = 48.18244161069647
$data[20,5] = 3.686666668227583
$data[20,6] = 0
$data[20,7] = 30.75472367724059
$data[20,8] = 10.64205923896326
$data[20,9] = 0
$data[20,10] = 0
$data[20,11] = 20.14841723888014

$data[21,0] = 18.65990287726171
$data[21,1] = 16.23693738332973
$data[21,2] = 6.909765081834202
$data[21,3] = 0
$data[21,4] = 48.01106705268977
$data[21,5] = 3.688260102177888
$data[21,6] = 0
$data[21,7] = 30.69108276837944
$data[21,8] = 10.64779673633268
$data[21,9] = 0
$data[21,10] = 0
$data[21,11] = 20.09031538993327

$data[22,0] = 17.98362360709296
$data[22,1] = 15.54052193610548
$data[22,2] = 6.922638945407401
$data[22,3] = 0
$data[22,4] = 47.37812686128344
$data[22,5] = 3.694515439409687
$data[22,6] = 0
$data[22,7] = 30.46102392413959
$data[22,8] = 10.67178379102664
$data[22,9] = 0
$data[22,10] = 0
$data[22,11] = 19.87743087647064

$data[23,0] = 17.24482915324912
$data[23,1] = 14.77083594839479
$data[23,2] = 6.938366872143018
$data[23,3] = 0
$data[23,4] = 46.73219304935753
$data[23,5] = 3.70173839683555
$data[23,6] = 0
$data[23,7] = 30.23675833783629
$data[23,8] = 10.70247452242763
$data[23,9] = 0
$data[23,10] = 0
$data[23,11] = 19.66400392183867

$ws.Range("B2:M25").Value = $data

